$wb = $excel.ActiveWorkbook

# --- Sheet rename: "Include from TRE-R13-Commune" -> "Include from TRE-R13-CommuneO"
$wsInclude = $wb.Worksheets.Item("Include from TRE-R13-Commune")
$wsInclude.Name = "Include from TRE-R13-CommuneO"

# --- Metadata sheet updates
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 2.0.1 -> 2.1.0
$wsMeta.Range("B4").Value = "2.1.0"

# Date: 2024-04-16T11:49:14+02:00 -> 2024-09-04T10:06:33+02:00
$wsMeta.Range("B9").Value = "2024-09-04T10:06:33+02:00"

# Contact: InteropSante (fhir@interopsante.org(WORK)) -> InteropSante (fhir@interopsante.org(work))
$wsMeta.Range("B12").Value = "InteropSanté (fhir@interopsante.org(work))"

# --- Include sheet updates (System URI row)
$wsInclude.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R13-CommuneOM/FHIR/TRE-R13-CommuneOM"
